$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: set NumberFormat to Text before assigning values that
# look numeric (contain dots/percent signs), then restore the original "Normal"
# style so the saved XML does not pick up a stray style index.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.617.84'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('E2').NumberFormat = 'General'
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.150.74'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E3').NumberFormat = 'General'
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E4').NumberFormat = 'General'
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '612.87'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('E5').NumberFormat = 'General'
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.08'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('E6').NumberFormat = 'General'
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E7').NumberFormat = 'General'
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.147.51'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('E8').NumberFormat = 'General'
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('E9').NumberFormat = 'General'
$ws.Range('E9').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('E10').NumberFormat = 'General'
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.51'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.93%  '
$ws.Range('E11').NumberFormat = 'General'
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.471'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.68%  '
$ws.Range('E12').NumberFormat = 'General'
$ws.Range('E12').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('E13').NumberFormat = 'General'
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.72'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.17%  '
$ws.Range('E14').NumberFormat = 'General'
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.673.50'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('E15').NumberFormat = 'General'
$ws.Range('E15').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('E16').NumberFormat = 'General'
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.562.39'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('E17').NumberFormat = 'General'
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.154.81'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('E18').NumberFormat = 'General'
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.93'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('E19').NumberFormat = 'General'
$ws.Range('E19').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('E20').NumberFormat = 'General'
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.68'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('E21').NumberFormat = 'General'
$ws.Range('E21').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('E22').NumberFormat = 'General'
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.92%  '
$ws.Range('E23').NumberFormat = 'General'
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.72'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('E24').NumberFormat = 'General'
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.00'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.42%  '
$ws.Range('E25').NumberFormat = 'General'
$ws.Range('E25').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E26').NumberFormat = 'General'
$ws.Range('E26').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.18%  '
$ws.Range('E27').NumberFormat = 'General'
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.54'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.31%  '
$ws.Range('E28').NumberFormat = 'General'
$ws.Range('E28').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.72%  '
$ws.Range('E29').NumberFormat = 'General'
$ws.Range('E29').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('E30').NumberFormat = 'General'
$ws.Range('E30').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.10%  '
$ws.Range('E31').NumberFormat = 'General'
$ws.Range('E31').Style = 'Normal'

$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Stacks'
$ws.Range('B32').NumberFormat = 'General'
$ws.Range('B32').Style = 'Normal'

$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C32').NumberFormat = 'General'
$ws.Range('C32').Style = 'Normal'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.73'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('E32').NumberFormat = 'General'
$ws.Range('E32').Style = 'Normal'

$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('B33').NumberFormat = 'General'
$ws.Range('B33').Style = 'Normal'

$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C33').NumberFormat = 'General'
$ws.Range('C33').Style = 'Normal'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E33').NumberFormat = 'General'
$ws.Range('E33').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.18%  '
$ws.Range('E34').NumberFormat = 'General'
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.15'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.26%  '
$ws.Range('E35').NumberFormat = 'General'
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0788'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.84%  '
$ws.Range('E36').NumberFormat = 'General'
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.01'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.69%  '
$ws.Range('E37').NumberFormat = 'General'
$ws.Range('E37').Style = 'Normal'

$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('B38').NumberFormat = 'General'
$ws.Range('B38').Style = 'Normal'

$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C38').NumberFormat = 'General'
$ws.Range('C38').Style = 'Normal'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.21'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.97%  '
$ws.Range('E38').NumberFormat = 'General'
$ws.Range('E38').Style = 'Normal'

$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'OKB'
$ws.Range('B39').NumberFormat = 'General'
$ws.Range('B39').Style = 'Normal'

$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C39').NumberFormat = 'General'
$ws.Range('C39').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '53.22'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('E39').NumberFormat = 'General'
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '458.75'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('E40').NumberFormat = 'General'
$ws.Range('E40').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('E41').NumberFormat = 'General'
$ws.Range('E41').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.48%  '
$ws.Range('E42').NumberFormat = 'General'
$ws.Range('E42').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('E43').NumberFormat = 'General'
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.846.52'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.96%  '
$ws.Range('E44').NumberFormat = 'General'
$ws.Range('E44').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.92%  '
$ws.Range('E45').NumberFormat = 'General'
$ws.Range('E45').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.60%  '
$ws.Range('E46').NumberFormat = 'General'
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.47'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +4.19%  '
$ws.Range('E47').NumberFormat = 'General'
$ws.Range('E47').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.26%  '
$ws.Range('E48').NumberFormat = 'General'
$ws.Range('E48').Style = 'Normal'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.05%  '
$ws.Range('E50').NumberFormat = 'General'
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.51'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.32%  '
$ws.Range('E51').NumberFormat = 'General'
$ws.Range('E51').Style = 'Normal'
